$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ptprz1"
$ws.Range("C2").Value = "L1cam"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.138766
$ws.Range("H2").Value = 0.416298
$ws.Range("I2").Value = 0.01356925767068476
$ws.Range("J2").Value = 0.01356925767068476
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.603177
$ws.Range("N2").Value = 19.809531
$ws.Range("O2").Value = 0.5135477412645301
$ws.Range("P2").Value = 0.5135477412645302
$ws.Range("Q2").Value = 0.9162964595819999
$ws.Range("R2").Value = 8.246668136238
$ws.Range("S2").Value = 0.006968461627416556
$ws.Range("T2").Value = 0.006968461627416558

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ptprz1"
$ws.Range("C3").Value = "L1cam"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.138766
$ws.Range("H3").Value = 0.416298
$ws.Range("I3").Value = 0.01356925767068476
$ws.Range("J3").Value = 0.01356925767068476
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3785896666666667
$ws.Range("N3").Value = 1.135769
$ws.Range("O3").Value = 0.02944398858046029
$ws.Range("P3").Value = 0.0294439885804603
$ws.Range("Q3").Value = 0.05253537368466667
$ws.Range("R3").Value = 0.472818363162
$ws.Range("S3").Value = 0.0003995330679009652
$ws.Range("T3").Value = 0.0003995330679009653

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ptprz1"
$ws.Range("C4").Value = "L1cam"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.138766
$ws.Range("H4").Value = 0.416298
$ws.Range("I4").Value = 0.01356925767068476
$ws.Range("J4").Value = 0.01356925767068476
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.441487333333333
$ws.Range("N4").Value = 10.324462
$ws.Range("O4").Value = 0.2676541983690312
$ws.Range("P4").Value = 0.2676541983690313
$ws.Range("Q4").Value = 0.4775614312973334
$ws.Range("R4").Value = 4.298052881676
$ws.Range("S4").Value = 0.003631868784309956
$ws.Range("T4").Value = 0.003631868784309958

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Ptprz1"
$ws.Range("C5").Value = "L1cam"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.138766
$ws.Range("H5").Value = 0.416298
$ws.Range("I5").Value = 0.01356925767068476
$ws.Range("J5").Value = 0.01356925767068476
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.434707333333333
$ws.Range("N5").Value = 7.304122
$ws.Range("O5").Value = 0.1893540717859783
$ws.Range("P5").Value = 0.1893540717859783
$ws.Range("Q5").Value = 0.3378545978173333
$ws.Range("R5").Value = 3.040691380356
$ws.Range("S5").Value = 0.002569394191057278
$ws.Range("T5").Value = 0.002569394191057279

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ptprz1"
$ws.Range("C6").Value = "L1cam"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.05416133333333333
$ws.Range("H6").Value = 0.162484
$ws.Range("I6").Value = 0.00529617548814441
$ws.Range("J6").Value = 0.005296175488144411
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.603177
$ws.Range("N6").Value = 19.809531
$ws.Range("O6").Value = 0.5135477412645301
$ws.Range("P6").Value = 0.5135477412645302
$ws.Range("Q6").Value = 0.357636870556
$ws.Range("R6").Value = 3.218731835004
$ws.Range("S6").Value = 0.002719838959277132
$ws.Range("T6").Value = 0.002719838959277133

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ptprz1"
$ws.Range("C7").Value = "L1cam"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.05416133333333333
$ws.Range("H7").Value = 0.162484
$ws.Range("I7").Value = 0.00529617548814441
$ws.Range("J7").Value = 0.005296175488144411
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.3785896666666667
$ws.Range("N7").Value = 1.135769
$ws.Range("O7").Value = 0.02944398858046029
$ws.Range("P7").Value = 0.0294439885804603
$ws.Range("Q7").Value = 0.02050492113288889
$ws.Range("R7").Value = 0.184544290196
$ws.Range("S7").Value = 0.0001559405305930377
$ws.Range("T7").Value = 0.0001559405305930378

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Ptprz1"
$ws.Range("C8").Value = "L1cam"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.05416133333333333
$ws.Range("H8").Value = 0.162484
$ws.Range("I8").Value = 0.00529617548814441
$ws.Range("J8").Value = 0.005296175488144411
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.441487333333333
$ws.Range("N8").Value = 10.324462
$ws.Range("O8").Value = 0.2676541983690312
$ws.Range("P8").Value = 0.2676541983690313
$ws.Range("Q8").Value = 0.1863955426231111
$ws.Range("R8").Value = 1.677559883608
$ws.Range("S8").Value = 0.001417543604701005
$ws.Range("T8").Value = 0.001417543604701005

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Ptprz1"
$ws.Range("C9").Value = "L1cam"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.05416133333333333
$ws.Range("H9").Value = 0.162484
$ws.Range("I9").Value = 0.00529617548814441
$ws.Range("J9").Value = 0.005296175488144411
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.434707333333333
$ws.Range("N9").Value = 7.304122
$ws.Range("O9").Value = 0.1893540717859783
$ws.Range("P9").Value = 0.1893540717859783
$ws.Range("Q9").Value = 0.1318669954497778
$ws.Range("R9").Value = 1.186802959048
$ws.Range("S9").Value = 0.001002852393573235
$ws.Range("T9").Value = 0.001002852393573236

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Ptprz1"
$ws.Range("C10").Value = "L1cam"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 10.002366
$ws.Range("H10").Value = 30.007098
$ws.Range("I10").Value = 0.9780831152479456
$ws.Range("J10").Value = 0.9780831152479456
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 6.603177
$ws.Range("N10").Value = 19.809531
$ws.Range("O10").Value = 0.5135477412645301
$ws.Range("P10").Value = 0.5135477412645302
$ws.Range("Q10").Value = 66.047393116782
$ws.Range("R10").Value = 594.4265380510379
$ws.Range("S10").Value = 0.5022923746045574
$ws.Range("T10").Value = 0.5022923746045576

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Ptprz1"
$ws.Range("C11").Value = "L1cam"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 10.002366
$ws.Range("H11").Value = 30.007098
$ws.Range("I11").Value = 0.9780831152479456
$ws.Range("J11").Value = 0.9780831152479456
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.3785896666666667
$ws.Range("N11").Value = 1.135769
$ws.Range("O11").Value = 0.02944398858046029
$ws.Range("P11").Value = 0.0294439885804603
$ws.Range("Q11").Value = 3.786792409818
$ws.Range("R11").Value = 34.081131688362
$ws.Range("S11").Value = 0.02879866807610154
$ws.Range("T11").Value = 0.02879866807610154

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Ptprz1"
$ws.Range("C12").Value = "L1cam"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 10.002366
$ws.Range("H12").Value = 30.007098
$ws.Range("I12").Value = 0.9780831152479456
$ws.Range("J12").Value = 0.9780831152479456
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.441487333333333
$ws.Range("N12").Value = 10.324462
$ws.Range("O12").Value = 0.2676541983690312
$ws.Range("P12").Value = 0.2676541983690313
$ws.Range("Q12").Value = 34.423015892364
$ws.Range("R12").Value = 309.807143031276
$ws.Range("S12").Value = 0.2617880521499736
$ws.Range("T12").Value = 0.2617880521499737

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Ptprz1"
$ws.Range("C13").Value = "L1cam"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 10.002366
$ws.Range("H13").Value = 30.007098
$ws.Range("I13").Value = 0.9780831152479456
$ws.Range("J13").Value = 0.9780831152479456
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.434707333333333
$ws.Range("N13").Value = 7.304122
$ws.Range("O13").Value = 0.1893540717859783
$ws.Range("P13").Value = 0.1893540717859783
$ws.Range("Q13").Value = 24.352833850884
$ws.Range("R13").Value = 219.175504657956
$ws.Range("S13").Value = 0.1852040204173128
$ws.Range("T13").Value = 0.1852040204173128

# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Ptprz1"
$ws.Range("C14").Value = "L1cam"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.03120566666666667
$ws.Range("H14").Value = 0.09361700000000001
$ws.Range("I14").Value = 0.003051451593225274
$ws.Range("J14").Value = 0.003051451593225274
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 6.603177
$ws.Range("N14").Value = 19.809531
$ws.Range("O14").Value = 0.5135477412645301
$ws.Range("P14").Value = 0.5135477412645302
$ws.Range("Q14").Value = 0.206056540403
$ws.Range("R14").Value = 1.854508863627
$ws.Range("S14").Value = 0.001567066073278891
$ws.Range("T14").Value = 0.001567066073278891

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Ptprz1"
$ws.Range("C15").Value = "L1cam"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.03120566666666667
$ws.Range("H15").Value = 0.09361700000000001
$ws.Range("I15").Value = 0.003051451593225274
$ws.Range("J15").Value = 0.003051451593225274
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.3785896666666667
$ws.Range("N15").Value = 1.135769
$ws.Range("O15").Value = 0.02944398858046029
$ws.Range("P15").Value = 0.0294439885804603
$ws.Range("Q15").Value = 0.01181414294144444
$ws.Range("R15").Value = 0.106327286473
$ws.Range("S15").Value = 0.00008984690586475233
$ws.Range("T15").Value = 0.00008984690586475233

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Ptprz1"
$ws.Range("C16").Value = "L1cam"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.03120566666666667
$ws.Range("H16").Value = 0.09361700000000001
$ws.Range("I16").Value = 0.003051451593225274
$ws.Range("J16").Value = 0.003051451593225274
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 3.441487333333333
$ws.Range("N16").Value = 10.324462
$ws.Range("O16").Value = 0.2676541983690312
$ws.Range("P16").Value = 0.2676541983690313
$ws.Range("Q16").Value = 0.1073939065615556
$ws.Range("R16").Value = 0.9665451590540001
$ws.Range("S16").Value = 0.0008167338300466137
$ws.Range("T16").Value = 0.0008167338300466139

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Ptprz1"
$ws.Range("C17").Value = "L1cam"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.03120566666666667
$ws.Range("H17").Value = 0.09361700000000001
$ws.Range("I17").Value = 0.003051451593225274
$ws.Range("J17").Value = 0.003051451593225274
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 2.434707333333333
$ws.Range("N17").Value = 7.304122
$ws.Range("O17").Value = 0.1893540717859783
$ws.Range("P17").Value = 0.1893540717859783
$ws.Range("Q17").Value = 0.07597666547488889
$ws.Range("R17").Value = 0.683789989274
$ws.Range("S17").Value = 0.0005778047840350163
$ws.Range("T17").Value = 0.0005778047840350164

